$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row for student "d06edc19" (Davis) - row 2
$ws.Rows(2).Delete()

# Fix the GPA for "342ab1a5" (Sanchez), now shifted up to row 2: 3.54 -> 3.82
$ws.Cells.Item(2, 5).Value = 3.82

# Fix the GPA for "3c346d6e" (Green), now shifted up to row 3: 3.7 -> 3.54
$ws.Cells.Item(3, 5).Value = 3.54
